# Add a new "2021" column (Y) to the worksheet, mirroring the existing
# year columns (D..X) for 2000..2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new year column.
$ws.Range("Y4").Value = 2021

# New data values for the 2021 column, one per data row.
$ws.Range("Y5").Value = 46.69
$ws.Range("Y6").Value = 52.52
$ws.Range("Y7").Value = 43.22
$ws.Range("Y8").Value = 51.31
$ws.Range("Y9").Value = 41.31
$ws.Range("Y10").Value = 52.43
$ws.Range("Y11").Value = 49.27
$ws.Range("Y12").Value = 31.68
$ws.Range("Y13").Value = 35.59
$ws.Range("Y14").Value = 55.28
$ws.Range("Y15").Value = 61.02
$ws.Range("Y16").Value = 48.72

# Match the formatting of the neighboring (X) column for each row so the
# new column inherits the same number formats/styles.
$ws.Range("Y4").NumberFormat = $ws.Range("X4").NumberFormat
for ($r = 5; $r -le 16; $r++) {
    $ws.Range("Y$r").NumberFormat = $ws.Range("X$r").NumberFormat
}

# Update the view: scroll so column B is the leftmost visible column and
# select cell AA15, matching the author's saved view state.
$ws.Range("AA15").Select()
$excel.ActiveWindow.ScrollColumn = 2
